# Wildfires sheet: add a "Main Table" summary on Home, and extend the
# Wildfires sheet with a wildfire-specific header row plus new Table 6-9
# sub-tables (Airports/Ports/Dams -> renumbered, plus a new Nuclear Power
# Plant table). Mirrors the "Modified the code to include script and
# created Py.ipynb is complete for Wildfires" commit.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Home sheet: new "Main Table" index pointing at the Wildfires table.
# (Shared strings for these new labels must be minted first, in this
# left-to-right / top-to-bottom order, so the shared string table lines
# up with the target file.)
# ------------------------------------------------------------------
$home = $wb.Worksheets.Item("Home")
$home.Range("A1").Value = "Main Table"
$home.Range("A3").Value = "Event_id "
$home.Range("B3").Value = "Episode"
$home.Range("C3").Value = "Event type"
$home.Range("D3").Value = "URL "

# ------------------------------------------------------------------
# Wildfires sheet: rework the main header row + renumber/extend the
# Table 6-9 sub-tables, finishing with a new Nuclear Power Plant table.
# ------------------------------------------------------------------
$wf = $wb.Worksheets.Item("Wildfires")

$wf.Range("C3").Value = "Countries"
$wf.Range("D3").Value = "Start_date/last_detected"
$wf.Range("E3").Value = "Duration"
$wf.Range("F3").Value = "People affected "
$wf.Range("G3").Value = "Burned area"
$wf.Range("H3").Value = "More info"

$wf.Range("A13").Value = "Table 4"
$wf.Range("A17").Value = "Table 5"
$wf.Range("A21").Value = "Table 6"
$wf.Range("A25").Value = "Table 7"
$wf.Range("A29").Value = "Table 8"

$wf.Range("A33").Value = "Table 9"
$wf.Range("B33").Value = "Nuclear Power Plant"

$wf.Range("C35").Value = "Name"
$wf.Range("D35").Value = "Country"
$wf.Range("E35").Value = "Reactors"
$wf.Range("F35").Value = "Distance"

# ------------------------------------------------------------------
# View-state updates recorded in the workbook when it was last saved.
# ------------------------------------------------------------------
$home.Range("E3").Select()

$eq = $wb.Worksheets.Item("Earthquake")
$eq.Activate()
$eq.Range("E9").Select()

$wf.Activate()
$excel.ActiveWindow.Zoom = 85
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$wf.Range("G5").Select()
